$wb = $excel.ActiveWorkbook

# --- Sheet 1: Overview ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = "ffff29a62014-0460-405d-88a9-e8168a4ad209.md"
$ws1.Range("A3").Value = "ffffff0b2dfd1d-7929-493c-a689-d8b26d1c65ab.md"
$ws1.Range("A4").Value = "796bf38a-be79-44d8-9419-f9c75750f5d6.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

# --- Sheet 2: zh-cn ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = "ffff29a62014-0460-405d-88a9-e8168a4ad209.md"
$ws2.Range("C2").Value = "a35ce3dd-932b-43fa-be8c-68e0f7f0c54d.57064ced9f1784210164ffd6a2d120012a081521.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-03-11 01:32:56"
$ws2.Range("E2").Value = "a35ce3dd-932b-43fa-be8c-68e0f7f0c54d.md"
$ws2.Range("F2").Value = "a35ce3dd-932b-43fa-be8c-68e0f7f0c54d.57064ced9f1784210164ffd6a2d120012a081521.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-03-11 01:33:43"

$ws2.Range("A3").Value = "ffffff0b2dfd1d-7929-493c-a689-d8b26d1c65ab.md"

$ws2.Range("A4").Value = "796bf38a-be79-44d8-9419-f9c75750f5d6.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "796bf38a-be79-44d8-9419-f9c75750f5d6.2289f68baeca0c7d5cd5b3828faa510e815b14b9.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-03-11 01:38:36"
$ws2.Range("E4").Value = "796bf38a-be79-44d8-9419-f9c75750f5d6.md"
$ws2.Range("F4").Value = "796bf38a-be79-44d8-9419-f9c75750f5d6.2289f68baeca0c7d5cd5b3828faa510e815b14b9.zh-cn.xlf"
$ws2.Range("G4").Value = "2016-03-11 01:37:41"

# --- Sheet 3: de-de ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = "ffff29a62014-0460-405d-88a9-e8168a4ad209.md"
$ws3.Range("C2").Value = "a35ce3dd-932b-43fa-be8c-68e0f7f0c54d.57064ced9f1784210164ffd6a2d120012a081521.de-de.xlf"
$ws3.Range("D2").Value = "2016-03-11 01:33:03"
$ws3.Range("E2").Value = "a35ce3dd-932b-43fa-be8c-68e0f7f0c54d.md"
$ws3.Range("F2").Value = "a35ce3dd-932b-43fa-be8c-68e0f7f0c54d.57064ced9f1784210164ffd6a2d120012a081521.de-de.xlf"
$ws3.Range("G2").Value = "2016-03-11 01:34:04"

$ws3.Range("A3").Value = "ffffff0b2dfd1d-7929-493c-a689-d8b26d1c65ab.md"

$ws3.Range("A4").Value = "796bf38a-be79-44d8-9419-f9c75750f5d6.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "796bf38a-be79-44d8-9419-f9c75750f5d6.2289f68baeca0c7d5cd5b3828faa510e815b14b9.de-de.xlf"
$ws3.Range("D4").Value = "2016-03-11 01:38:43"
$ws3.Range("E4").Value = "796bf38a-be79-44d8-9419-f9c75750f5d6.md"
$ws3.Range("F4").Value = "796bf38a-be79-44d8-9419-f9c75750f5d6.2289f68baeca0c7d5cd5b3828faa510e815b14b9.de-de.xlf"
$ws3.Range("G4").Value = "2016-03-11 01:38:02"
